$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.219.40"
$ws.Range("D3").Value = "1.828.73"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "236.51"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").Value = "0.6083"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.07105"
$ws.Range("E8").Value = "  -4.45%  "
$ws.Range("D9").Value = "0.2814"
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("D10").Value = "23.74"
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").Value = "0.07663"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "1.822.68"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "4.828"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "0.00001008"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "0.6329"
$ws.Range("E15").Value = "  -6.32%  "
$ws.Range("D16").Value = "2.072.95"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "79.19"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "5.879"
$ws.Range("E18").Value = "  -5.86%  "
$ws.Range("D19").Value = "29.219.86"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").Value = "228.09"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.021"
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "154.91"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("D26").Value = "8.065"
$ws.Range("E26").Value = "  -4.99%  "
$ws.Range("D27").Value = "0.1306"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").Value = "1.488"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").Value = "0.06445"
$ws.Range("E30").Value = "  -6.70%  "
$ws.Range("D31").Value = "1.453"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "3.825"
$ws.Range("E32").Value = "  -5.48%  "
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("D34").Value = "1.127"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "1.752"
$ws.Range("E36").Value = "  -6.64%  "
$ws.Range("D37").Value = "2.547"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "2.754"
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").Value = "1.220.78"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").Value = "6.562"
$ws.Range("E41").Value = "  -3.57%  "
$ws.Range("D42").Value = "0.9313"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "101.13"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "1.982.88"
$ws.Range("D46").Value = "63.07"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.615"
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "8.604"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05527"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1076"
$ws.Range("E51").Value = "  -5.54%  "
